$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "modelo"
$ws.Range("B1").Value = "marca"
$ws.Range("C1").Value = "cor"
$ws.Range("D1").Value = "placa"
$ws.Range("E1").Value = "ano"
$ws.Range("F1").Value = "observacoes"
$ws.Range("G1").Value = "valor de compra"
$ws.Range("H1").Value = "status"
$ws.Range("I1").Value = "Data de Cadastro"
$ws.Range("J1").Value = "Valor diaria"

# --- Row 2 (existing record, updated) ---
$ws.Range("A2").Value = "Yamaha i8"
$ws.Range("B2").Value = "Yamaha"
$ws.Range("C2").Value = "preto"
$ws.Range("D2").Value = "JDSA214"
$ws.Range("E2").Value = 2025
$ws.Range("F2").Value = "teste"
$ws.Range("G2").Value = 10000
$ws.Range("H2").Value = "Ativo"
$ws.Range("I2").Value = "2024-05-19 09:57:11"
$ws.Range("J2").Value = 50

# --- Row 3 (new record) ---
$ws.Range("A3").Value = "Teste"
$ws.Range("B3").Value = "teste"
$ws.Range("C3").Value = "Azul Metálico"
$ws.Range("D3").Value = "tiue822"
$ws.Range("E3").Value = 2020
$ws.Range("F3").Value = "teste"
$ws.Range("G3").Value = 14000
$ws.Range("H3").Value = "Ativo"
$ws.Range("I3").Value = "2024-05-19 10:23:20"
$ws.Range("J3").Value = 50

# --- Autofilter range now covers A1:J3 ---
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A1:J3").AutoFilter()

# --- Update the hidden _FilterDatabase defined name to match the new range ---
$fdb = $wb.Names.Item("ConfigMotos!_FilterDatabase")
$fdb.RefersTo = "='ConfigMotos'!`$A`$1:`$J`$3"
